$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

$data = @(
    @{A="bCTzf322"; B=23102880; C="ugefrpi66"; D="cg`$!5T9K"; E="MR"; F="usvEzSSc"; G="HIDF"; H="Candidate"},
    @{A="sDFaW246"; B=23102879; C="pksbmwu45"; D="N&V#5va9"; E="MR"; F="ktVXlKhU"; G="kqkw"; H="Candidate"},
    @{A="UYOEX176"; B=23102878; C="msqeuul45"; D="A7eD#c2%"; E="MR"; F="UHKTMaOD"; G="IPYe"; H="Candidate"},
    @{A="uSbMW814"; B=23102877; C="kuecmli75"; D="ufA3Y&`$7"; E="MR"; F="bYhDwXMQ"; G="PsUu"; H="Candidate"},
    @{A="kKpqN200"; B=23102876; C="rrtuhjq30"; D="X!2Vdc4&"; E="MR"; F="gykPksXn"; G="WLUs"; H="Candidate"},
    @{A="WolfQ745"; B=23102875; C="lvtofbk89"; D="wx8#D&E6"; E="MR"; F="aiKWdGCS"; G="NldZ"; H="Candidate"},
    @{A="qYrnS116"; B=23102874; C="xiuyakn78"; D="bGV%8`$6z"; E="MR"; F="hEsFesKm"; G="cqtT"; H="Candidate"},
    @{A="tVBKZ849"; B=23102873; C="llfmulp17"; D="Vg&b5Y4#"; E="MR"; F="LHTkPAXT"; G="UHJv"; H="Candidate"},
    @{A="Anzdf811"; B=23102872; C="uriuzib17"; D="WZh54b&!"; E="MR"; F="AxDWCbWz"; G="dFtr"; H="Candidate"},
    @{A="CnRyK282"; B=23102871; C="qibyeoj54"; D="dq#4T%M3"; E="MR"; F="otSRVTrr"; G="HKdd"; H="Candidate"},
    @{A="uvwln554"; B=23102870; C="xhexxwf68"; D="n7Vp!#D8"; E="MR"; F="JbCaczdK"; G="ufaZ"; H="Candidate"},
    @{A="crlLW260"; B=23102869; C="iiqpeer25"; D="dV%bM49!"; E="MR"; F="BpUDvXZO"; G="cDov"; H="Candidate"},
    @{A="OHxqD566"; B=23102868; C="iayvmsn63"; D="Au%8Ct&9"; E="MR"; F="piNvMwbA"; G="kXnQ"; H="Candidate"}
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
}

# Delete row 15 which is no longer needed
$ws.Rows.Item(15).Delete()
